$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "계획 확장"
$ws.Range("F4").Value = "계획 확장 모드에서 사용하는 도구들입니다."
$ws.Range("F5").Value = "계획 보이기 전환"
$ws.Range("F6").Value = "색상 선택"
$ws.Range("F7").Value = "수정"
$ws.Range("F8").Value = "덮어쓰기 금지"
$ws.Range("F9").Value = "계획 (문)"
$ws.Range("F10").Value = "계획 (바닥)"
$ws.Range("F11").Value = "계획 (구조물)"
$ws.Range("F12").Value = "계획 (벽)"
$ws.Range("F13").Value = "형상 변환"
$ws.Range("F14").Value = "작업 1"
$ws.Range("F15").Value = "작업 2"
$ws.Range("F16").Value = "되돌리기"
$ws.Range("F17").Value = "다시 실행"
$ws.Range("F21").Value = "문 형태의 계획을 배치합니다.\n\n우클릭으로 형태를 선택할 수 있습니다.\n\n{0} 키를 눌러 색상을 선택할 수 있습니다.\n\n{1} 키를 눌러 크기를 조절할 수 있습니다."
$ws.Range("F23").Value = "바닥 형태의 계획을 배치합니다.\n\n우클릭으로 형태를 선택할 수 있습니다.\n\n{0} 키를 눌러 색상을 선택할 수 있습니다.\n\n{1} 키를 눌러 크기를 조절할 수 있습니다."
$ws.Range("F25").Value = "구조물 형태의 계획을 배치합니다.\n\n우클릭으로 형태를 선택할 수 있습니다.\n\n{0} 키를 눌러 색상을 선택할 수 있습니다.\n\n{1} 키를 눌러 크기를 조절할 수 있습니다."
$ws.Range("F27").Value = "벽 형태의 계획을 배치합니다.\n\n우클릭으로 형태를 선택할 수 있습니다.\n\n{0} 키를 눌러 색상을 선택할 수 있습니다.\n\n{1} 키를 눌러 크기를 조절할 수 있습니다."
$ws.Range("F29").Value = "이미 존재하는 계획을 채색합니다.\n\n우클릭으로 형태를 선택할 수 있습니다.\n\n{0} 키를 눌러 색상을 선택할 수 있습니다.\n\n{1} 키를 눌러 크기를 조절할 수 있습니다."
$ws.Range("F31").Value = "계획을 제거합니다."
$ws.Range("F33").Value = "선택한 영역의 계획을 복사합니다."
$ws.Range("F35").Value = "선택한 영역의 계획을 잘라내 복사합니다."
$ws.Range("F37").Value = "저장된 계획을 붙여넣습니다.\n\n키를 사용해 회전하거나 뒤집을 수 있습니다."
$ws.Range("F39").Value = "마지막으로 생성한 계획 작업을 취소합니다."
$ws.Range("F41").Value = "마지막으로 취소한 계획 작업을 다시 실행합니다."
$ws.Range("F43").Value = "파일에서 계획을 불러옵니다.\n\n좌클릭으로 계획 목록을 열 수 있습니다.\n\n우클릭으로 마지막으로 불러온 계획을 빠르게 선택할 수 있습니다."
$ws.Range("F45").Value = "마지막으로 복사한 계획을 파일에 저장합니다."
$ws.Range("F47").Value = "모든 계획의 가시성을 전환합니다.\n\n(0)을 누른 상태로 각 계획을 선택하면 개별로 전환할 수 있습니다."
$ws.Range("F49").Value = "좌클릭으로 불투명도를 설정합니다.\n\n우클릭으로 텍스처 세트를 변경합니다.\n\n(0)을 누른 상태로 각 계획을 선택하면 개별로 전환할 수 있습니다."
$ws.Range("F50").Value = "색상 계획을 선택합니다."
$ws.Range("F51").Value = "기존 계획에서 색상을 가져옵니다."
$ws.Range("F52").Value = "계획을 찾을 수 없습니다."
$ws.Range("F53").Value = "계획이 복사되었습니다."
$ws.Range("F54").Value = "계획을 삭제합니다."
$ws.Range("F55").Value = "계획을 {0}(으)로 저장합니다."
$ws.Range("F56").Value = "{0} {1} 계획을 삭제합니다."
$ws.Range("F57").Value = "문"
$ws.Range("F58").Value = "바닥"
$ws.Range("F59").Value = "구조물"
$ws.Range("F60").Value = "벽"
$ws.Range("F61").Value = "모드"
$ws.Range("F62").Value = "복사"
$ws.Range("F63").Value = "잘라내기"
$ws.Range("F64").Value = "건너뛰기"
$ws.Range("F65").Value = "교체"
$ws.Range("F66").Value = "잘라내기로 전환"
$ws.Range("F67").Value = "되돌리기-다시 실행"
$ws.Range("F69").Value = "되돌리기-다시 실행 기능을 사용하고 버튼이 표시될지 여부를 결정합니다."
$ws.Range("F72").Value = "활성화하면 계획 잘라내기 버튼이 표시됩니다. 비활성화하면 보이지 않습니다."
$ws.Range("F74").Value = "활성화하면 계획 보이기 버튼이 표시됩니다. 비활성화하면 보이지 않습니다."
$ws.Range("F75").Value = "보기 설정 버튼 표시"
$ws.Range("F76").Value = "활성화하면 보기 설정 버튼이 표시됩니다. 비활성화하면 보이지 않습니다."
$ws.Range("F78").Value = "활성화하면 건설을 완료한 후에도 계획이 유지됩니다. 비활성화하면 제거됩니다."
$ws.Range("F79").Value = "컨트롤 키로 색상 띄우기"
$ws.Range("F80").Value = "활성화하면 컨트롤 키를 눌러 색상 목록을 표시할 수 있습니다."
$ws.Range("F81").Value = "교체 대신 건너뛰기를 기본값으로 사용"
$ws.Range("F82").Value = "활성화하면 기본적으로 계획을 건너뜁니다. 비활성화하면 교체를 기본값으로 사용합니다."
$ws.Range("F83").Value = "계획 표시"
$ws.Range("F84").Value = "계획의 표시 여부를 결정합니다."
$ws.Range("F113").Value = "{0} 개수"
$ws.Range("F116").Value = "넓이"
$ws.Range("F118").Value = "없음"
$ws.Range("F119").Value = "북"
$ws.Range("F120").Value = "북동"
$ws.Range("F121").Value = "동"
$ws.Range("F122").Value = "남동"
$ws.Range("F123").Value = "남"
$ws.Range("F124").Value = "남서"
$ws.Range("F125").Value = "서"
$ws.Range("F126").Value = "북서"
$ws.Range("F129").Value = "대각"
$ws.Range("F130").Value = "북서 대각"
$ws.Range("F131").Value = "북동 대각"

$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F121").Select()
